$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32-42 down to 33-43.
$ws.Rows(32).Insert()

# Populate the new row 32 with a new weekly price observation (matches the
# other rows for this market/product: same A-C, E-K, Q, R columns).
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(32, 3).Value = "Ñuble"
$ws.Cells.Item(32, 4).Value = 45006
$ws.Cells.Item(32, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(32, 5).Value = 16
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100101
$ws.Cells.Item(32, 8).Value = "Berries"
$ws.Cells.Item(32, 9).Value = 100101001
$ws.Cells.Item(32, 10).Value = "Arándano (blue)"
$ws.Cells.Item(32, 11).Value = "Sin especificar"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 60
$ws.Cells.Item(32, 14).Value = 4000
$ws.Cells.Item(32, 15).Value = 4000
$ws.Cells.Item(32, 16).Value = 4000
$ws.Cells.Item(32, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(32, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(32, 19).Value = 2000
$ws.Cells.Item(32, 20).Value = 2
